$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Shift rows 12-15 down to 13-16 to make room for a new "Jurisdiction" row ---
# (Copy() on this host only overwrites non-blank destination cells, so clear
# the destination first to make sure blank source cells really blank it out.)
$ws.Range("A16:B16").ClearContents()
$ws.Range("A15:B15").Copy($ws.Range("A16"))

$ws.Range("A15:B15").ClearContents()
$ws.Range("A14:B14").Copy($ws.Range("A15"))

$ws.Range("A14:B14").ClearContents()
$ws.Range("A13:B13").Copy($ws.Range("A14"))

$ws.Range("A13:B13").ClearContents()
$ws.Range("A12:B12").Copy($ws.Range("A13"))

# Row 12 becomes the new "Jurisdiction" property (value left blank), reusing
# the formatting that was already in place on that row.
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# --- Field value updates ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
